# Add two new columns ("I0" / "IF") to the sheet, matching the style of
# the existing header row, and fill in their data values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (same bold/bordered style as the rest of row 1).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the formatting (style) from the existing "IP" header cell (H1)
# onto the two new header cells so they reuse the same cell style.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

# New data values for row 2 (unstyled, like the other data cells).
$ws.Range("I2").Value = 5
$ws.Range("J2").Value = 7
